$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('I2').Value = 6690
$ws.Range('I3').Value = 6983
$ws.Range('F4').Value = 1872
$ws.Range('I4').Value = 1595
$ws.Range('I6').Value = 8067
$ws.Range('F7').Value = 24061
$ws.Range('I7').Value = 23985

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range('I4').Value = 41
$ws.Range('I7').Value = 285

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range('I2').Value = 80
$ws.Range('I3').Value = 63
$ws.Range('I6').Value = 108
$ws.Range('I7').Value = 276

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range('I2').Value = 121
$ws.Range('I7').Value = 426

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('I3').Value = 342
$ws.Range('I7').Value = 915

$ws = $wb.Worksheets.Item('New City')
$ws.Range('I4').Value = 22
$ws.Range('I6').Value = 168
$ws.Range('I7').Value = 557

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('I2').Value = 187
$ws.Range('I8').Value = 1436
$ws.Range('I11').Value = 365
$ws.Range('I12').Value = 61
$ws.Range('I19').Value = 675
$ws.Range('I20').Value = 593
$ws.Range('I23').Value = 233
$ws.Range('I29').Value = 1446
$ws.Range('I33').Value = 1067
$ws.Range('I36').Value = 329
$ws.Range('I40').Value = 42
$ws.Range('I47').Value = 173
$ws.Range('I49').Value = 160
$ws.Range('I50').Value = 123
$ws.Range('I53').Value = 262
$ws.Range('I57').Value = 94
$ws.Range('I59').Value = 39
$ws.Range('F63').Value = 162
$ws.Range('I63').Value = 72
$ws.Range('I65').Value = 557
$ws.Range('I67').Value = 915
$ws.Range('I73').Value = 218
$ws.Range('I75').Value = 77
$ws.Range('I77').Value = 144
$ws.Range('I79').Value = 683
$ws.Range('I85').Value = 1076
$ws.Range('I89').Value = 285
$ws.Range('I92').Value = 72
$ws.Range('I94').Value = 245
$ws.Range('I96').Value = 276
$ws.Range('I99').Value = 426
$ws.Range('F101').Value = 24061
$ws.Range('I101').Value = 23985

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('I6').Value = 345
$ws.Range('I7').Value = 1067

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range('I2').Value = 29
$ws.Range('I7').Value = 160

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('I3').Value = 501
$ws.Range('I4').Value = 76
$ws.Range('I7').Value = 1446

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range('I3').Value = 198
$ws.Range('I7').Value = 675

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('I6').Value = 281
$ws.Range('I7').Value = 1076

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range('I2').Value = 63
$ws.Range('I7').Value = 233

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range('I2').Value = 198
$ws.Range('I7').Value = 683

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range('I3').Value = 169
$ws.Range('I7').Value = 593

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range('I3').Value = 109
$ws.Range('I7').Value = 329

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range('I2').Value = 47
$ws.Range('I7').Value = 245

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range('I3').Value = 51
$ws.Range('I7').Value = 173

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range('I2').Value = 37
$ws.Range('I7').Value = 123

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range('I2').Value = 144
$ws.Range('I3').Value = 81
$ws.Range('I6').Value = 97
$ws.Range('I7').Value = 365

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range('I6').Value = 58
$ws.Range('I7').Value = 218

$ws = $wb.Worksheets.Item('Montclare')
$ws.Range('I3').Value = 4
$ws.Range('I7').Value = 39

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range('I2').Value = 69
$ws.Range('I7').Value = 187

$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Range('I2').Value = 23
$ws.Range('I7').Value = 72

$ws = $wb.Worksheets.Item('Austin')
$ws.Range('I2').Value = 427
$ws.Range('I6').Value = 461
$ws.Range('I7').Value = 1436

$ws = $wb.Worksheets.Item('Pullman')
$ws.Range('I2').Value = 24
$ws.Range('I7').Value = 77

$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Range('I2').Value = 35
$ws.Range('I7').Value = 94

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range('I6').Value = 123
$ws.Range('I7').Value = 262

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range('I2').Value = 47
$ws.Range('I6').Value = 35
$ws.Range('I7').Value = 144

$ws = $wb.Worksheets.Item('Hegewisch')
$ws.Range('I6').Value = 7
$ws.Range('I7').Value = 42

$ws = $wb.Worksheets.Item('Beverly')
$ws.Range('I6').Value = 33
$ws.Range('I7').Value = 61

"Updated 106 cells"
